# Apply trade #6 close update to the live trading results workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Summary sheet - update aggregate metrics
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1500.14   # Current Capital
$summary.Range("B4").Value = 0.14      # Total P&L $
$summary.Range("B5").Value = 0.47      # Total P&L %
$summary.Range("B6").Value = 6         # Total Trades
$summary.Range("B7").Value = 3         # Winning Trades
$summary.Range("B9").Value = 50        # Win Rate %

# ---------------------------------------------------------------------------
# 2. Strategy Status sheet - update MarketMaking strategy row (row 6)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C6").Value = 100.14     # Capital
$status.Range("D6").Value = 6          # Trades
$status.Range("E6").Value = 0.14       # P&L $
$status.Range("F6").Value = 0.14       # P&L %
$status.Range("G6").Value = 50         # Win Rate %

# ---------------------------------------------------------------------------
# 3. Append new trade #6 row to "All Trades" and "MarketMaking" sheets
# ---------------------------------------------------------------------------
# Note: the Date column value is prefixed with an apostrophe so Excel keeps
# it as literal text ("2026-02-17") instead of auto-converting it to a date
# serial number, matching the rest of the column.
$newRow = @(6, "'2026-02-17", "23:52:26", "MarketMaking", "UP", 0.91, 0.95, "CLOSED", 4.3956, 0.04, 100.14, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.13)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 7
    for ($i = 0; $i -lt $newRow.Count; $i++) {
        $col = $i + 1
        $ws.Cells.Item($row, $col).Value = $newRow[$i]
    }
}
